$d = $word.ActiveDocument

# 1) Paragraph about seminar content: trim "những ... đóng góp" and
#    rephrase "và điều này sẽ giúp" -> "đồng thời cũng sẽ giúp"
$d.Content.Find.Execute(
    "chia sẻ những ý kiến đóng góp của mình và điều này sẽ giúp",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "chia sẻ ý kiến của mình đồng thời cũng sẽ giúp", 2)

# 2) "Vui lòng phản hồi chúng tôi" -> "Vui lòng xác nhận tham gia"
$d.Content.Find.Execute(
    "Vui lòng phản hồi chúng tôi bằng cách gửi đơn đăng ký trước ngày",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Vui lòng xác nhận tham gia bằng cách gửi đơn đăng ký trước ngày", 2)

# 3) Button text: "Gửi thông tin cá nhân" -> "Gửi thông tin của tôi"
$d.Content.Find.Execute(
    "Gửi thông tin cá nhân",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Gửi thông tin của tôi", 2)
